# "planilha com dados é ida direitin"
# The data block (A1:A13 and B2:F13 — i.e. the populated data range minus the
# still-"general" B1:F1 header cells and the untouched G column) gets
# re-aligned from its previous center/general look to right-aligned, and a
# handful of cells in the bottom data row (row 13) / row 11 get their values
# corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Alignment: switch the data block to right-aligned ---
# xlRight = -4152
$ws.Range("A1:A13").HorizontalAlignment = -4152
$ws.Range("B2:F13").HorizontalAlignment = -4152

# --- Value corrections ---
$ws.Range("D11").Value = 0
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 1
$ws.Range("F13").Value = 20
